$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Add new row 2 data (order chosen to match shared-string insertion order:
# OBT, Verifying the invitations, Y, RCC010)
$ws.Range("B2").Value = "OBT"
$ws.Range("C2").Value = "Verifying the invitations"
$ws.Range("D2").Value = "Y"
$ws.Range("A2").Value = "RCC010"

# Apply a thin border (matching the existing header border) to the new row,
# without the header's bold font / yellow fill.
$rng = $ws.Range("A2:E2")
$rng.Borders.ColorIndex = 1
$rng.Borders.Weight = 2

# Update the view: move the selection (this also resets any scrolled
# topLeftCell back to the default visible area)
$ws.Range("B7").Select()
